$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the style of the existing header row (bold, bordered, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 98
    $ws.Cells.Item($r, 31).Value = 64
    $ws.Cells.Item($r, 32).Value = 0
}
